$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9467159249493364
$ws.Range("D2").Value = 0.9634230801028263
$ws.Range("C3").Value = 0.9010723949085274
$ws.Range("D3").Value = 0.930922589278271
$ws.Range("C4").Value = 0.8617743837261826
$ws.Range("D4").Value = 0.9021660098948185
$ws.Range("C5").Value = 0.8280555113172309
$ws.Range("D5").Value = 0.8768806214540861
$ws.Range("C6").Value = 0.7986857967825759
$ws.Range("D6").Value = 0.8540877782506606
$ws.Range("C7").Value = 0.7729072301579705
$ws.Range("D7").Value = 0.8334593397229496
$ws.Range("C8").Value = 0.7499555415707542
$ws.Range("D8").Value = 0.8145821739549657
$ws.Range("C9").Value = 0.7300861805253056
$ws.Range("D9").Value = 0.7975162488513825
$ws.Range("C10").Value = 0.7124612553818958
$ws.Range("D10").Value = 0.7822123504499183
$ws.Range("C11").Value = 0.6974442412558072
$ws.Range("D11").Value = 0.7695344664708723
$ws.Range("C12").Value = 0.6836666960791763
$ws.Range("D12").Value = 0.7578079566845245
$ws.Range("C13").Value = 0.6719133269603114
$ws.Range("D13").Value = 0.7471380062646644
$ws.Range("C14").Value = 0.6611930416409698
$ws.Range("D14").Value = 0.7375688464355227
$ws.Range("C15").Value = 0.6517088735068035
$ws.Range("D15").Value = 0.7283894148336741
$ws.Range("C16").Value = 0.6430971470711536
$ws.Range("D16").Value = 0.7197894106815088
$ws.Range("C17").Value = 0.6353267967027401
$ws.Range("D17").Value = 0.7129042250671698
$ws.Range("C18").Value = 0.6287898037223315
$ws.Range("D18").Value = 0.706956707521074
$ws.Range("C19").Value = 0.6226453260041828
$ws.Range("D19").Value = 0.7015195646640769
$ws.Range("B20").Value = 0.655916913233514
$ws.Range("C20").Value = 0.6161611768505176
$ws.Range("D20").Value = 0.6956655103999354
$ws.Range("C21").Value = 0.6109818967028349
$ws.Range("D21").Value = 0.6902462582193687
$ws.Range("C22").Value = 0.6063284426610934
$ws.Range("D22").Value = 0.6860383211551058
$ws.Range("C23").Value = 0.5979747007948456
$ws.Range("D23").Value = 0.6779859484732571
$ws.Range("C24").Value = 0.5942497247614262
$ws.Range("D24").Value = 0.6746203347395816
$ws.Range("C25").Value = 0.591139310158724
$ws.Range("D25").Value = 0.6716157895944079
$ws.Range("C26").Value = 0.588437376174694
$ws.Range("D26").Value = 0.6689523463241377
$ws.Range("C27").Value = 0.5832881056601664
$ws.Range("D27").Value = 0.6643503306523221
$ws.Range("C28").Value = 0.5809138161902065
$ws.Range("D28").Value = 0.6623251535486798
$ws.Range("B29").Value = 0.620715358082193
$ws.Range("C29").Value = 0.5786878947844895
$ws.Range("D29").Value = 0.6607209799551258
$ws.Range("C30").Value = 0.5766628893854987
$ws.Range("D30").Value = 0.6594255728190823
$ws.Range("C31").Value = 0.5715746908032616
$ws.Range("D31").Value = 0.6562175181442607
$ws.Range("C32").Value = 0.5701671103672132
$ws.Range("D32").Value = 0.6552016142221485
$ws.Range("C33").Value = 0.5682538588507728
$ws.Range("D33").Value = 0.6528742385214078
$ws.Range("B34").Value = 0.6107630141623641
$ws.Range("C34").Value = 0.5672420376869004
$ws.Range("D34").Value = 0.6520640631895687
$ws.Range("C35").Value = 0.5664299902701584
$ws.Range("D35").Value = 0.6514886357990374
$ws.Range("C36").Value = 0.5630377363742297
$ws.Range("D36").Value = 0.6475791046713434
$ws.Range("C37").Value = 0.5618616123647372
$ws.Range("D37").Value = 0.6463995858923486
$ws.Range("C38").Value = 0.5608907860445151
$ws.Range("D38").Value = 0.6452604981626902
$ws.Range("C39").Value = 0.5583736327625133
$ws.Range("D39").Value = 0.643940655808804
